$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before writing, so that numeric-looking
# strings like "4.716" or "28.720.19" are preserved as text and not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.720.19'
$ws.Range("E2").Value = '  +6.90%  '
$ws.Range("D3").Value = '1.809.49'
$ws.Range("E3").Value = '  +4.67%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '251.54'
$ws.Range("E5").Value = '  +4.05%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '0.4969'
$ws.Range("E7").Value = '  +1.45%  '
$ws.Range("D8").Value = '0.2801'
$ws.Range("E8").Value = '  +7.75%  '
$ws.Range("E9").Value = '  +2.75%  '
$ws.Range("D10").Value = '1.802.11'
$ws.Range("D11").Value = '16.78'
$ws.Range("E11").Value = '  +4.71%  '
$ws.Range("D12").Value = '0.07109'
$ws.Range("E12").Value = '  +2.82%  '
$ws.Range("D13").Value = '0.6500'
$ws.Range("E13").Value = '  +6.56%  '
$ws.Range("D14").Value = '4.716'
$ws.Range("E14").Value = '  +4.93%  '
$ws.Range("D15").Value = '81.94'
$ws.Range("E15").Value = '  +5.94%  '
$ws.Range("D16").Value = '28.703.92'
$ws.Range("D17").Value = '0.9994'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '0.000007357'
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").Value = '0.9988'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '12.27'
$ws.Range("E20").Value = '  +7.13%  '
$ws.Range("D21").Value = '2.039.08'
$ws.Range("E21").Value = '  +4.17%  '
$ws.Range("D22").Value = '4.617'
$ws.Range("E22").Value = '  +4.08%  '
$ws.Range("D23").Value = '8.904'
$ws.Range("E23").Value = '  +3.90%  '
$ws.Range("D24").Value = '5.315'
$ws.Range("E24").Value = '  +3.64%  '
$ws.Range("D25").Value = '142.97'
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("D26").Value = '16.02'
$ws.Range("E26").Value = '  +4.46%  '
$ws.Range("D27").Value = '1.887'
$ws.Range("E27").Value = '  +5.94%  '
$ws.Range("D28").Value = '112.09'
$ws.Range("E28").Value = '  +5.37%  '
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = '4.191'
$ws.Range("E30").Value = '  +5.91%  '
$ws.Range("D31").Value = '0.08369'
$ws.Range("E31").Value = '  +4.65%  '
$ws.Range("D32").Value = '3.843'
$ws.Range("E32").Value = '  +4.22%  '
$ws.Range("D33").Value = '0.04984'
$ws.Range("E33").Value = '  +9.90%  '
$ws.Range("D34").Value = '1.094'
$ws.Range("E34").Value = '  +8.38%  '
$ws.Range("D35").Value = '0.6752'
$ws.Range("E35").Value = '  +8.05%  '
$ws.Range("D36").Value = '2.666'
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("D37").Value = '0.9638'
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("D38").Value = '2.669'
$ws.Range("E38").Value = '  +9.31%  '
$ws.Range("D39").Value = '2.149'
$ws.Range("D40").Value = '0.01594'
$ws.Range("E40").Value = '  +5.98%  '
$ws.Range("D41").Value = '5.938'
$ws.Range("E41").Value = '  +4.78%  '
$ws.Range("D42").Value = '0.9994'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4123'
$ws.Range("E43").Value = '  +6.66%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '101.06'
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("D45").Value = '7.225'
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("D46").Value = '0.1225'
$ws.Range("E46").Value = '  +5.44%  '
$ws.Range("D47").Value = '0.05492'
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("D48").Value = '8.206'
$ws.Range("E48").Value = '  +3.06%  '
$ws.Range("D49").Value = '31.44'
$ws.Range("E49").Value = '  +4.00%  '
$ws.Range("D50").Value = '1.303'
$ws.Range("E50").Value = '  +5.18%  '
$ws.Range("E51").Value = '  +6.62%  '

# Restore the original (default/Normal) style on column D now that the text values
# are locked in, so no stray number-format style lingers on these cells.
$ws.Range("D2:D51").Style = "Normal"
